# Clear the converted "retailer_id" and "additional_emails" example rows
# from the template (rows 2-17), remove their mailto: hyperlinks, and
# select rows 2-19 as the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the sample retailer_id values (column A) and the sample
# additional_emails values (column C) for rows 2 through 17.
$ws.Range("A2:A17").ClearContents()
$ws.Range("C2:C17").ClearContents()

# Remove the mailto: hyperlinks that were attached to the now-cleared
# email cells.
$ws.Hyperlinks.Delete()

# Reflect the selection state saved with the workbook (rows 2-19 selected).
$ws.Rows("2:19").Select() | Out-Null
